# Updates the "Estado de Cuenta" worker detail table:
#  - Rows 16-18 now show the new worker "DOLIS MANUEL ESPINOZA JIMENEZ" (doc 1044391428)
#    for periods 1909, 1910, 1911 with Valor Mora 33125 and Salario Basico 877803.
#  - Rows 19-37 keep "MATILDE ISABEL HERRERA MIRANDA" (doc 34967348) but now list her
#    periods in ascending order (2205..2311) with an updated Salario Basico of 1423500.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New first block (3 rows): DOLIS MANUEL ESPINOZA JIMENEZ
$dolisPeriods = @("1909", "1910", "1911")
for ($i = 0; $i -lt $dolisPeriods.Length; $i++) {
    $r = 16 + $i
    $ws.Cells.Item($r, 3).Value = "1044391428"
    $ws.Cells.Item($r, 4).Value = "DOLIS MANUEL ESPINOZA JIMENEZ"
    $ws.Cells.Item($r, 5).Value = $dolisPeriods[$i]
    $ws.Cells.Item($r, 6).Value = 33125
    $ws.Cells.Item($r, 7).Value = 877803
}

# Second block (19 rows): MATILDE ISABEL HERRERA MIRANDA, periods now ascending
$matildePeriods = @("2205", "2206", "2207", "2208", "2209", "2210", "2211", "2212", "2301", "2302", "2303", "2304", "2305", "2306", "2307", "2308", "2309", "2310", "2311")
for ($i = 0; $i -lt $matildePeriods.Length; $i++) {
    $r = 19 + $i
    $ws.Cells.Item($r, 3).Value = "34967348"
    $ws.Cells.Item($r, 4).Value = "MATILDE ISABEL HERRERA MIRANDA"
    $ws.Cells.Item($r, 5).Value = $matildePeriods[$i]
    $ws.Cells.Item($r, 6).Value = 40268
    $ws.Cells.Item($r, 7).Value = 1423500
}
